# Natmi following Dr Hou advice
# The LR-pair result table for Efna1-Epha1 now includes "ECs" (endothelial cells)
# as an additional sending/target cluster alongside the existing "FAPs" and "sCs",
# so the 3x3 sending/target cluster combinations grow the table from 6 to 9 data rows
# (rows 2-10) and every NATMI specificity/expression statistic is recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Efna1-Epha1)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 21.78783233333333
$ws.Range("H2").Value = 65.363497
$ws.Range("I2").Value = 0.9024488799587679
$ws.Range("J2").Value = 0.9024488799587679
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.160945333333332
$ws.Range("N2").Value = 24.482836
$ws.Range("O2").Value = 0.491005088714322
$ws.Range("P2").Value = 0.491005088714322
$ws.Range("Q2").Value = 177.8093086041657
$ws.Range("R2").Value = 1600.283777437492
$ws.Range("S2").Value = 0.4431069923642953
$ws.Range("T2").Value = 0.4431069923642953

# Row 3: ECs -> FAPs (Efna1-Epha1)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 21.78783233333333
$ws.Range("H3").Value = 65.363497
$ws.Range("I3").Value = 0.9024488799587679
$ws.Range("J3").Value = 0.9024488799587679
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.477462333333333
$ws.Range("N3").Value = 7.432386999999999
$ws.Range("O3").Value = 0.1490570715865667
$ws.Range("P3").Value = 0.1490570715865668
$ws.Range("Q3").Value = 53.97853393081543
$ws.Range("R3").Value = 485.8068053773389
$ws.Range("S3").Value = 0.134516387303231
$ws.Range("T3").Value = 0.134516387303231

# Row 4: ECs -> sCs (Efna1-Epha1)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 21.78783233333333
$ws.Range("H4").Value = 65.363497
$ws.Range("I4").Value = 0.9024488799587679
$ws.Range("J4").Value = 0.9024488799587679
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.982489999999999
$ws.Range("N4").Value = 17.94747
$ws.Range("O4").Value = 0.3599378396991113
$ws.Range("P4").Value = 0.3599378396991114
$ws.Range("Q4").Value = 130.3454890558433
$ws.Range("R4").Value = 1173.10940150259
$ws.Range("S4").Value = 0.3248255002912415
$ws.Range("T4").Value = 0.3248255002912416

# Row 5: FAPs -> ECs (Efna1-Epha1)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.359006333333333
$ws.Range("H5").Value = 4.077019
$ws.Range("I5").Value = 0.05628984676448105
$ws.Range("J5").Value = 0.05628984676448104
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.160945333333332
$ws.Range("N5").Value = 24.482836
$ws.Range("O5").Value = 0.491005088714322
$ws.Range("P5").Value = 0.491005088714322
$ws.Range("Q5").Value = 11.09077639398711
$ws.Range("R5").Value = 99.81698754588399
$ws.Range("S5").Value = 0.02763860120430961
$ws.Range("T5").Value = 0.0276386012043096

# Row 6: FAPs -> FAPs (Efna1-Epha1)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.359006333333333
$ws.Range("H6").Value = 4.077019
$ws.Range("I6").Value = 0.05628984676448105
$ws.Range("J6").Value = 0.05628984676448104
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.477462333333333
$ws.Range("N6").Value = 7.432386999999999
$ws.Range("O6").Value = 0.1490570715865667
$ws.Range("P6").Value = 0.1490570715865668
$ws.Range("Q6").Value = 3.366887001594777
$ws.Range("R6").Value = 30.301983014353
$ws.Range("S6").Value = 0.008390399718770123
$ws.Range("T6").Value = 0.008390399718770123

# Row 7: FAPs -> sCs (Efna1-Epha1)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.359006333333333
$ws.Range("H7").Value = 4.077019
$ws.Range("I7").Value = 0.05628984676448105
$ws.Range("J7").Value = 0.05628984676448104
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.982489999999999
$ws.Range("N7").Value = 17.94747
$ws.Range("O7").Value = 0.3599378396991113
$ws.Range("P7").Value = 0.3599378396991114
$ws.Range("Q7").Value = 8.130241799103333
$ws.Range("R7").Value = 73.17217619192999
$ws.Range("S7").Value = 0.02026084584140132
$ws.Range("T7").Value = 0.02026084584140132

# Row 8: sCs -> ECs (Efna1-Epha1)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9961713333333334
$ws.Range("H8").Value = 2.988514
$ws.Range("I8").Value = 0.04126127327675106
$ws.Range("J8").Value = 0.04126127327675105
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.160945333333332
$ws.Range("N8").Value = 24.482836
$ws.Range("O8").Value = 0.491005088714322
$ws.Range("P8").Value = 0.491005088714322
$ws.Range("Q8").Value = 8.129699793967111
$ws.Range("R8").Value = 73.16729814570401
$ws.Range("S8").Value = 0.02025949514571704
$ws.Range("T8").Value = 0.02025949514571703

# Row 9: sCs -> FAPs (Efna1-Epha1)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9961713333333334
$ws.Range("H9").Value = 2.988514
$ws.Range("I9").Value = 0.04126127327675106
$ws.Range("J9").Value = 0.04126127327675105
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.477462333333333
$ws.Range("N9").Value = 7.432386999999999
$ws.Range("O9").Value = 0.1490570715865667
$ws.Range("P9").Value = 0.1490570715865668
$ws.Range("Q9").Value = 2.467976955879778
$ws.Range("R9").Value = 22.211792602918
$ws.Range("S9").Value = 0.006150284564565575
$ws.Range("T9").Value = 0.006150284564565575

# Row 10: sCs -> sCs (Efna1-Epha1)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9961713333333334
$ws.Range("H10").Value = 2.988514
$ws.Range("I10").Value = 0.04126127327675106
$ws.Range("J10").Value = 0.04126127327675105
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.982489999999999
$ws.Range("N10").Value = 17.94747
$ws.Range("O10").Value = 0.3599378396991113
$ws.Range("P10").Value = 0.3599378396991114
$ws.Range("Q10").Value = 5.959585039953333
$ws.Range("R10").Value = 53.63626535958
$ws.Range("S10").Value = 0.01485149356646845
$ws.Range("T10").Value = 0.01485149356646845

